# Update the cached "datetimeFigureOut" date field text shown on every
# Date Placeholder (slide master + all slide layouts) from 06/03/2015
# to 03/09/2015.

$p = $ppt.ActivePresentation

$oldDate = "06/03/2015"
$newDate = "03/09/2015"

# ppPlaceholderDate
$ppPlaceholderDate = 16

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePh = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePh = $true
            }
        } catch {
            $isDatePh = $false
        }
        if ($isDatePh -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide Master
$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

# Every slide layout (CustomLayout) under the master
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholders $layout.Shapes
}
